# Auto-generated edit script applying the diff changes
# (gh-pages data refresh: updated "想去人数" counts and two cover image URLs)
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 101
$ws.Range("F4").Value = 8084
$ws.Range("F6").Value = 82
$ws.Range("F7").Value = 80
$ws.Range("F8").Value = 7010
$ws.Range("F9").Value = 7010
$ws.Range("F10").Value = 1132
$ws.Range("F11").Value = 530
$ws.Range("F12").Value = 482
$ws.Range("F14").Value = 697
$ws.Range("F16").Value = 305
$ws.Range("F21").Value = 11473
$ws.Range("F22").Value = 114
$ws.Range("F23").Value = 2218
$ws.Range("F25").Value = 3061
$ws.Range("F26").Value = 51
$ws.Range("F27").Value = 44
$ws.Range("F28").Value = 2642
$ws.Range("F29").Value = 99
$ws.Range("F31").Value = 274
$ws.Range("F32").Value = 43
$ws.Range("F34").Value = 2346
$ws.Range("F36").Value = 1595
$ws.Range("F37").Value = 72
$ws.Range("I37").Value = "//i0.hdslb.com/bfs/openplatform/202407/c5hRL8Nq1721038527962.jpeg"
$ws.Range("F38").Value = 91
$ws.Range("F39").Value = 5763
$ws.Range("F40").Value = 1771
$ws.Range("F41").Value = 1242
$ws.Range("F42").Value = 823
$ws.Range("F43").Value = 157
$ws.Range("F46").Value = 1101
$ws.Range("F47").Value = 1063
$ws.Range("F48").Value = 1509
$ws.Range("F49").Value = 95
$ws.Range("F50").Value = 1126

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 249
$ws.Range("F11").Value = 209
$ws.Range("F20").Value = 62

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 210
$ws.Range("F3").Value = 344

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 101
$ws.Range("F4").Value = 210
$ws.Range("F5").Value = 344
$ws.Range("F8").Value = 8084
$ws.Range("F9").Value = 82
$ws.Range("F11").Value = 80
$ws.Range("F12").Value = 7010
$ws.Range("F13").Value = 1132
$ws.Range("F14").Value = 530
$ws.Range("F15").Value = 482
$ws.Range("F16").Value = 697
$ws.Range("F18").Value = 305
$ws.Range("F22").Value = 209
$ws.Range("F24").Value = 11473
$ws.Range("F25").Value = 114
$ws.Range("F26").Value = 2218
$ws.Range("F27").Value = 2218
$ws.Range("F28").Value = 3061
$ws.Range("F29").Value = 2642
$ws.Range("F31").Value = 274
$ws.Range("F32").Value = 43
$ws.Range("F34").Value = 2346
$ws.Range("F36").Value = 1595
$ws.Range("F37").Value = 72
$ws.Range("I37").Value = "//i0.hdslb.com/bfs/openplatform/202407/c5hRL8Nq1721038527962.jpeg"
$ws.Range("F38").Value = 91
$ws.Range("F39").Value = 5763
$ws.Range("F40").Value = 62
$ws.Range("F41").Value = 1771
$ws.Range("F43").Value = 1242
$ws.Range("F44").Value = 823
$ws.Range("F45").Value = 157
$ws.Range("F47").Value = 1101
$ws.Range("F48").Value = 1063
$ws.Range("F49").Value = 1509
$ws.Range("F50").Value = 95
$ws.Range("F51").Value = 1126

